# Applies the "feedback 2.0" hour-log entries to rows 25 and 26,
# and updates the selection/scroll position to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("A25").Value = "Added a subview to the pinView"
$ws.Range("B25").Value = (Get-Date -Year 2018 -Month 1 -Day 3)
$ws.Range("C25").Value = 0.5
$ws.Range("D25").Value = "I made popupAdress subclassed to pinView. Now the popup will move with the pinView"

# Row 26
$ws.Range("A26").Value = "Research pinView "
$ws.Range("B26").Value = (Get-Date -Year 2018 -Month 1 -Day 3)
$ws.Range("C26").Value = 1.5
$ws.Range("D26").Value = "I did research to change the position of the popupAdress while subclassed to the pinView but I couldnt find the solution"

# Scroll / selection state
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("D25").Select()
